# chart_result.xlsx edit: "modify last-tested version and notes"
#
# - Rename the worksheet from "Multiple Queries" to "Chart"
# - Update the "(Last tested with: ...)" note to reflect the new
#   ReportServer version (4.0.0-6053 -> 4.1.0-6064)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Multiple Queries" -> "Chart"
$ws.Name = "Chart"

# Keep the embedded chart's series references pointing at the renamed sheet
# (mirrors what Excel does automatically when a referenced sheet is renamed).
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Chart!`$C`$7,Chart!`$A`$8:`$A`$13,Chart!`$C`$8:`$C`$13,1)"

# Update the last-tested-with note (row 19, column A)
$ws.Range("A19").Value = "(Last tested with: ReportServer 4.1.0-6064) "
